# This script fixes a batch of mistranslated / placeholder labels in the
# "16.9.1" SDG indicator worksheet (birth registration of children under 5).
#
# The header cells below were inadvertently left with English "(*)"-style
# stub text (e.g. "Sex", "Male", "Female", "Areas", "Education", lower-case
# section labels, etc.) or with Kyrgyz/Russian column headers that didn't
# match the corrected wording used elsewhere in the sheet ("Эркек"/"Аял" vs.
# the plural "Эркектер"/"Аялдар" used for the group header, "Возраст (в
# месяцах)" vs. "По возрасту (в месяцах)", etc.). This replaces each of
# them with the corrected text, leaving every other cell (values, styles,
# merged layout) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "By sex" group header / sub-headers (row 6-8) ---
$ws.Range("C6").Value  = "By sex"
$ws.Range("A7").Value  = "Эркектер"
$ws.Range("B7").Value  = "Мужчины"
$ws.Range("C7").Value  = "Men"
$ws.Range("A8").Value  = "Аялдар"
$ws.Range("B8").Value  = "Женщины"
$ws.Range("C8").Value  = "Woman"

# --- "By territory" group header (row 12) ---
$ws.Range("C12").Value = "By territory"

# --- "By age (in month)" group header (row 22) ---
$ws.Range("A22").Value = "Жаш курагы боюнча (айларда)"
$ws.Range("B22").Value = "По возрасту (в месяцах)"
$ws.Range("C22").Value = "By age (in month)"

# --- "Education of mother" group + its English sub-items (rows 28-33) ---
$ws.Range("C28").Value = "Education of mother"
$ws.Range("C29").Value = "Preschool or not /primary"
$ws.Range("C30").Value = "Basic general"
$ws.Range("C31").Value = "Average total"
$ws.Range("C32").Value = "Vocational primary /secondary"
$ws.Range("C33").Value = "Higher"
